# Re-running everything with 5,000 runs (was 10,000) -- update the resulting
# percent-difference values in the summary table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 64
    "D2" = 37
    "E2" = 1
    "G2" = 30
    "C4" = 143
    "E4" = 3
    "G4" = 23
    "H4" = 24
    "I4" = 30
    "D5" = 36
    "E5" = 6
    "E6" = -15
    "F6" = -25
    "G6" = -70
    "H6" = -2
    "D7" = 50
    "G7" = 41
    "I7" = 52
    "H9" = -20
    "D10" = 12
    "E10" = 7
    "G10" = 10
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
